$wb = $excel.ActiveWorkbook

# --- Sheet: ALC (sheet index 1) ---
$ws = $wb.Worksheets.Item("ALC")

$ws.Range("H17").Value = 2050.24
$ws.Range("I17").Value = 1000
$ws.Range("J17").Value = 2141.5652
$ws.Range("K17").Value = 3000
$ws.Range("L17").Value = 6424.6956
$ws.Range("M17").Value = -2832
$ws.Range("N17").Value = -6760.6956

$ws.Range("H64").Value = 7593.125
$ws.Range("I64").Value = 4582
$ws.Range("J64").Value = 9399.799999999999
$ws.Range("K64").Value = 4582
$ws.Range("L64").Value = 9399.799999999999
$ws.Range("M64").Value = -4334
$ws.Range("N64").Value = -9895.799999999999

$ws.Range("H67").Value = 7593.125
$ws.Range("I67").Value = 4582
$ws.Range("J67").Value = 9399.799999999999
$ws.Range("K67").Value = 4582
$ws.Range("L67").Value = 9399.799999999999
$ws.Range("M67").Value = -3724
$ws.Range("N67").Value = -11115.8

$ws.Range("H113").Value = 9651
$ws.Range("I113").Value = 9651
$ws.Range("J113").Value = 0
$ws.Range("K113").Value = 9651
$ws.Range("L113").Value = 0
$ws.Range("M113").Value = -6397
$ws.Range("N113").Value = ""

$ws.Range("H138").Value = 1810.08
$ws.Range("I138").Value = 972.0625
$ws.Range("J138").Value = 3299.889
$ws.Range("K138").Value = 2916.1875
$ws.Range("L138").Value = 9899.667000000001
$ws.Range("M138").Value = 2223.8125
$ws.Range("N138").Value = -20179.667


# --- Sheet: ARM (sheet index 2) ---
$ws = $wb.Worksheets.Item("ARM")

$ws.Range("H32").Value = 2506363.5
$ws.Range("I32").Value = 3192.3547
$ws.Range("J32").Value = 11128398
$ws.Range("K32").Value = 3192.3547
$ws.Range("L32").Value = 11128398
$ws.Range("M32").Value = -2905.3547
$ws.Range("N32").Value = -11128972

$ws.Range("H43").Value = 7540085.5
$ws.Range("I43").Value = 10020114
$ws.Range("K43").Value = 10020114
$ws.Range("M43").Value = -10019801

$ws.Range("H63").Value = 6036.5
$ws.Range("I63").Value = 1600
$ws.Range("J63").Value = 10473
$ws.Range("K63").Value = 1600
$ws.Range("L63").Value = 10473
$ws.Range("M63").Value = -914
$ws.Range("N63").Value = -11845

$ws.Range("H66").Value = 6036.5
$ws.Range("I66").Value = 1600
$ws.Range("J66").Value = 10473
$ws.Range("K66").Value = 8000
$ws.Range("L66").Value = 52365
$ws.Range("M66").Value = -4568
$ws.Range("N66").Value = -59229

$ws.Range("H97").Value = 2303.2
$ws.Range("I97").Value = 1323.5
$ws.Range("K97").Value = 1323.5
$ws.Range("M97").Value = -827.5

$ws.Range("H102").Value = 4777.6665
$ws.Range("I102").Value = 2166.5
$ws.Range("K102").Value = 2166.5
$ws.Range("M102").Value = -544.5


# --- Sheet: BSM (sheet index 3) ---
$ws = $wb.Worksheets.Item("BSM")

$ws.Range("H2").Value = 59900
$ws.Range("J2").Value = 59900
$ws.Range("L2").Value = 59900
$ws.Range("N2").Value = -60126

$ws.Range("H20").Value = 3513.7144
$ws.Range("I20").Value = 3916.818
$ws.Range("J20").Value = 2035.6666
$ws.Range("K20").Value = 3916.818
$ws.Range("L20").Value = 2035.6666
$ws.Range("M20").Value = -3669.818
$ws.Range("N20").Value = -2529.6666

$ws.Range("H50").Value = 0
$ws.Range("J50").Value = 0
$ws.Range("L50").Value = 0
$ws.Range("N50").Value = ""

$ws.Range("H99").Value = 2331.8333
$ws.Range("I99").Value = 2331.8333
$ws.Range("K99").Value = 2331.8333
$ws.Range("M99").Value = -833.8332999999998

$ws.Range("H102").Value = 10555.667
$ws.Range("I102").Value = 10555.667
$ws.Range("K102").Value = 10555.667
$ws.Range("M102").Value = -7310.666999999999


# --- Sheet: CRP (sheet index 4) ---
$ws = $wb.Worksheets.Item("CRP")

$ws.Range("H3").Value = 2500950
$ws.Range("I3").Value = 5000450
$ws.Range("J3").Value = 1450
$ws.Range("K3").Value = 5000450
$ws.Range("L3").Value = 1450
$ws.Range("M3").Value = -5000337
$ws.Range("N3").Value = -1676

$ws.Range("H32").Value = 1252.8334
$ws.Range("I32").Value = 1252.8334
$ws.Range("K32").Value = 1252.8334
$ws.Range("M32").Value = -936.8334

$ws.Range("H43").Value = 26750
$ws.Range("J43").Value = 26750
$ws.Range("L43").Value = 26750
$ws.Range("N43").Value = -27118

$ws.Range("H58").Value = 6175.625
$ws.Range("I58").Value = 1441
$ws.Range("J58").Value = 6852
$ws.Range("K58").Value = 1441
$ws.Range("L58").Value = 6852
$ws.Range("M58").Value = -1238
$ws.Range("N58").Value = -7258

$ws.Range("H62").Value = 2889
$ws.Range("I62").Value = 2889
$ws.Range("K62").Value = 2889
$ws.Range("M62").Value = -2265

$ws.Range("H64").Value = 50000
$ws.Range("J64").Value = 50000
$ws.Range("L64").Value = 50000
$ws.Range("N64").Value = -50496

$ws.Range("H65").Value = 2889
$ws.Range("I65").Value = 2889
$ws.Range("K65").Value = 14445
$ws.Range("M65").Value = -11325

$ws.Range("H67").Value = 50000
$ws.Range("J67").Value = 50000
$ws.Range("L67").Value = 50000
$ws.Range("N67").Value = -51716

$ws.Range("H81").Value = 44832.668
$ws.Range("J81").Value = 44832.668
$ws.Range("L81").Value = 44832.668
$ws.Range("N81").Value = -46828.668

$ws.Range("H82").Value = 65000
$ws.Range("J82").Value = 65000
$ws.Range("L82").Value = 65000
$ws.Range("N82").Value = -65722

$ws.Range("H84").Value = 44832.668
$ws.Range("J84").Value = 44832.668
$ws.Range("L84").Value = 134498.004
$ws.Range("N84").Value = -144482.004

$ws.Range("H85").Value = 65000
$ws.Range("J85").Value = 65000
$ws.Range("L85").Value = 65000
$ws.Range("N85").Value = -67496

$ws.Range("H101").Value = 26750
$ws.Range("J101").Value = 26750
$ws.Range("L101").Value = 26750
$ws.Range("N101").Value = -33240

$ws.Range("H136").Value = 6175.625
$ws.Range("I136").Value = 1441
$ws.Range("J136").Value = 6852
$ws.Range("K136").Value = 4323
$ws.Range("L136").Value = 20556
$ws.Range("M136").Value = -1773
$ws.Range("N136").Value = -25656


# --- Sheet: CUL (sheet index 5) ---
$ws = $wb.Worksheets.Item("CUL")

$ws.Range("H21").Value = 69
$ws.Range("J21").Value = 0
$ws.Range("L21").Value = 0
$ws.Range("N21").Value = ""

$ws.Range("H44").Value = 825
$ws.Range("J44").Value = 1080
$ws.Range("L44").Value = 3240
$ws.Range("N44").Value = -4036

$ws.Range("H80").Value = 4588.346
$ws.Range("I80").Value = 4549.85
$ws.Range("K80").Value = 13649.55
$ws.Range("M80").Value = -12713.55

$ws.Range("H83").Value = 4588.346
$ws.Range("I83").Value = 4549.85
$ws.Range("K83").Value = 40948.65
$ws.Range("M83").Value = -36268.65

$ws.Range("H97").Value = 3022.6
$ws.Range("I97").Value = 1782.5
$ws.Range("K97").Value = 5347.5
$ws.Range("M97").Value = -4851.5

$ws.Range("H108").Value = 425
$ws.Range("I108").Value = 425
$ws.Range("J108").Value = 0
$ws.Range("K108").Value = 1275
$ws.Range("L108").Value = 0
$ws.Range("M108").Value = 1605
$ws.Range("N108").Value = ""

$ws.Range("H137").Value = 3966.3333
$ws.Range("J137").Value = 4449.5
$ws.Range("L137").Value = 13348.5
$ws.Range("N137").Value = -23548.5


# --- Sheet: GSM (sheet index 6) ---
$ws = $wb.Worksheets.Item("GSM")

$ws.Range("H70").Value = 8876.375
$ws.Range("I70").Value = 4601.2
$ws.Range("J70").Value = 16001.667
$ws.Range("K70").Value = 4601.2
$ws.Range("L70").Value = 16001.667
$ws.Range("M70").Value = -4331.2
$ws.Range("N70").Value = -16541.667

$ws.Range("H73").Value = 8876.375
$ws.Range("I73").Value = 4601.2
$ws.Range("J73").Value = 16001.667
$ws.Range("K73").Value = 4601.2
$ws.Range("L73").Value = 16001.667
$ws.Range("M73").Value = -3665.2
$ws.Range("N73").Value = -17873.667

$ws.Range("H109").Value = 44484
$ws.Range("J109").Value = 44484
$ws.Range("L109").Value = 44484
$ws.Range("N109").Value = -46564

$ws.Range("H113").Value = 6721.6313
$ws.Range("I113").Value = 4337.636
$ws.Range("K113").Value = 4337.636
$ws.Range("M113").Value = -2167.636

$ws.Range("H122").Value = 836189.7
$ws.Range("I122").Value = 1002429.6
$ws.Range("K122").Value = 3007288.8
$ws.Range("M122").Value = -3004838.8

$ws.Range("H126").Value = 6250
$ws.Range("I126").Value = 5000
$ws.Range("J126").Value = 7500
$ws.Range("K126").Value = 15000
$ws.Range("L126").Value = 22500
$ws.Range("M126").Value = -12530
$ws.Range("N126").Value = -27440


# --- Sheet: LTW (sheet index 7) ---
$ws = $wb.Worksheets.Item("LTW")

$ws.Range("H82").Value = 3892.7144
$ws.Range("I82").Value = 2265.5
$ws.Range("J82").Value = 5113.125
$ws.Range("K82").Value = 2265.5
$ws.Range("L82").Value = 5113.125
$ws.Range("M82").Value = -1904.5
$ws.Range("N82").Value = -5835.125

$ws.Range("H85").Value = 3892.7144
$ws.Range("I85").Value = 2265.5
$ws.Range("J85").Value = 5113.125
$ws.Range("K85").Value = 2265.5
$ws.Range("L85").Value = 5113.125
$ws.Range("M85").Value = -1017.5
$ws.Range("N85").Value = -7609.125

$ws.Range("H108").Value = 11569
$ws.Range("I108").Value = 11569
$ws.Range("K108").Value = 11569
$ws.Range("M108").Value = -7729

$ws.Range("H122").Value = 4898.3335
$ws.Range("I122").Value = 4885.625
$ws.Range("J122").Value = 5000
$ws.Range("K122").Value = 14656.875
$ws.Range("L122").Value = 15000
$ws.Range("M122").Value = -12206.875
$ws.Range("N122").Value = -19900


# --- Sheet: WVR (sheet index 8) ---
$ws = $wb.Worksheets.Item("WVR")

$ws.Range("H82").Value = 0
$ws.Range("J82").Value = 0
$ws.Range("L82").Value = 0
$ws.Range("N82").Value = ""

$ws.Range("H85").Value = 0
$ws.Range("J85").Value = 0
$ws.Range("L85").Value = 0
$ws.Range("N85").Value = ""

$ws.Range("H107").Value = 488.15384
$ws.Range("I107").Value = 376.9091
$ws.Range("K107").Value = 1130.7273
$ws.Range("M107").Value = 789.2727

